$wb = $excel.ActiveWorkbook

# --- Change 1: page27_table0, cells C7:F8 ---
# "as giving at a true 31 March and 2023 fair view of the entity's financial position performance for the financial year ended on that date, and"
# -> "giving a true and fair view of the entity's financial position as at 31 March 2023 performance for the financial year ended on that date, and"
$ws1 = $wb.Worksheets.Item("page27_table0")
$newText1 = "giving a true and fair view of the entity's financial position as at 31 March 2023 performance for the financial year ended on that date, and"
$ws1.Range("C7:F8").Value = $newText1

# --- Change 2: page13_table0, cells A26:E27 ---
# "Trade receivables are recognised initially at the amount of consideration that is unconditional unless they contain significant financing components, when they are recognised at fair value. The company holds the trade"
# -> "components, Trade they receivables company are recognised initially at the amount of consideration that is unconditional unless they contain significant financing when are recognised at fair value. The holds the trade"
$ws2 = $wb.Worksheets.Item("page13_table0")
$newText2 = "components, Trade they receivables company are recognised initially at the amount of consideration that is unconditional unless they contain significant financing when are recognised at fair value. The holds the trade"
$ws2.Range("A26:E27").Value = $newText2

# --- Change 3: page13_table0, cells A28:E29 ---
# "receivables with the objective of collecting the contractual cash flows and therefore them subsequently at amortised components, cost using they the effective interest method. company"
# -> "components, they company receivables with the objective of collecting the contractual cash flows and therefore them subsequently at amortised cost using the effective interest method."
$newText3 = "components, they company receivables with the objective of collecting the contractual cash flows and therefore them subsequently at amortised cost using the effective interest method."
$ws2.Range("A28:E29").Value = $newText3
